$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.030.33'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.609.16'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.10'
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.607.38'
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("E10").Value = '  -2.45%  '
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("E13").Value = '  -3.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.43'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.082.25'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("E16").Value = '  -3.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.969.55'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.607.97'
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '364.07'
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.97'
$ws.Range("E20").Value = '  -3.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.37'
$ws.Range("E21").Value = '  -5.03%  '
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("E23").Value = '  -1.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.97'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '66.39'
$ws.Range("E26").Value = '  -2.75%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.747.09'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '579.11'
$ws.Range("E29").Value = '  -2.78%  '
$ws.Range("E30").Value = '  -3.67%  '
$ws.Range("E31").Value = '  -4.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.73'
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("E33").Value = '  -2.17%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -8.42%  '
$ws.Range("E36").Value = '  -3.58%  '
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.31'
$ws.Range("E38").Value = '  +1.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.96'
$ws.Range("E39").Value = '  -3.20%  '
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.24'
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("E42").Value = '  -2.78%  '
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.37'
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '155.63'
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("E48").Value = '  -3.45%  '
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.63'
$ws.Range("E50").Value = '  +3.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.617'
$ws.Range("E51").Value = '  -2.22%  '
